$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume) are stored as plain text in the source
# workbook (e.g. "63.00", "1.00"), so force Text format before writing the
# new values to avoid Excel auto-converting numeric-looking strings to
# numbers (which would drop formatting like trailing zeros).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '26.631.48'
$ws.Range("E2").Value = '  -1.48%  '
$ws.Range("D3").Value = '1.591.72'
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '210.77'
$ws.Range("E5").Value = '  -1.71%  '
$ws.Range("E6").Value = '  -1.36%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("E8").Value = '  -1.97%  '
$ws.Range("E9").Value = '  -1.15%  '
$ws.Range("D10").Value = '19.63'
$ws.Range("E10").Value = '  -3.02%  '
$ws.Range("E11").Value = '  -1.34%  '
$ws.Range("D12").Value = '1.814.31'
$ws.Range("E12").Value = '  -1.73%  '
$ws.Range("D13").Value = '1.589.78'
$ws.Range("E13").Value = '  -1.94%  '
$ws.Range("E14").Value = '  -2.34%  '
$ws.Range("E15").Value = '  -3.20%  '
$ws.Range("D16").Value = '64.73'
$ws.Range("E16").Value = '  +0.59%  '
$ws.Range("D17").Value = '26.628.47'
$ws.Range("E17").Value = '  -1.39%  '
$ws.Range("D18").Value = '0.0₃0727'
$ws.Range("E18").Value = '  -1.91%  '
$ws.Range("D19").Value = '208.71'
$ws.Range("E19").Value = '  -3.31%  '
$ws.Range("E20").Value = '  -0.01%  '
$ws.Range("D21").Value = '6.73'
$ws.Range("E21").Value = '  -2.29%  '
$ws.Range("D22").Value = '4.25'
$ws.Range("E22").Value = '  -2.57%  '
$ws.Range("D23").Value = '2.37'
$ws.Range("E23").Value = '  -1.51%  '
$ws.Range("D24").Value = '8.88'
$ws.Range("E24").Value = '  -1.49%  '
$ws.Range("D25").Value = '146.82'
$ws.Range("E25").Value = '  -0.26%  '
$ws.Range("E26").Value = '  +0.04%  '
$ws.Range("D27").Value = '7.26'
$ws.Range("E27").Value = '  -0.13%  '
$ws.Range("E28").Value = '  -2.67%  '
$ws.Range("E29").Value = '  -1.57%  '
$ws.Range("D30").Value = '0.0509'
$ws.Range("E30").Value = '  +1.17%  '
$ws.Range("E31").Value = '  -1.90%  '
$ws.Range("D32").Value = '3.23'
$ws.Range("E32").Value = '  -3.42%  '
$ws.Range("E33").Value = '  +22.99%  '
$ws.Range("D34").Value = '2.91'
$ws.Range("E34").Value = '  -2.38%  '
$ws.Range("D35").Value = '1.315.52'
$ws.Range("E35").Value = '  -1.45%  '
$ws.Range("D36").Value = '2.42'
$ws.Range("E36").Value = '  -1.19%  '
$ws.Range("E37").Value = '  -4.24%  '
$ws.Range("E38").Value = '  -2.20%  '
$ws.Range("D39").Value = '0.831'
$ws.Range("E39").Value = '  -1.83%  '
$ws.Range("E40").Value = '  +0.06%  '
$ws.Range("D41").Value = '0.789'
$ws.Range("E41").Value = '  -1.39%  '
$ws.Range("E42").Value = '  +2.99%  '
$ws.Range("E43").Value = '  -2.66%  '
$ws.Range("D44").Value = '63.00'
$ws.Range("E44").Value = '  -2.64%  '
$ws.Range("D45").Value = '1.726.99'
$ws.Range("E45").Value = '  -1.57%  '
$ws.Range("D46").Value = '89.91'
$ws.Range("E46").Value = '  -0.41%  '
$ws.Range("E47").Value = '  +0.08%  '
$ws.Range("D48").Value = '0.829'
$ws.Range("E48").Value = '  -1.04%  '
$ws.Range("D49").Value = '0.0509'
$ws.Range("E49").Value = '  -0.53%  '
$ws.Range("D50").Value = '0.0981'
$ws.Range("E50").Value = '  -1.49%  '
$ws.Range("D51").Value = '7.53'
$ws.Range("E51").Value = '  -0.38%  '
